# LOQ4228.xlsx content update
# - Adds new "Objetivos" (PT) text, new "Docentes responsaveis" value
# - Shifts several rows' A-column labels down and re-wires the summary/
#   program/evaluation rows, adding a new Bibliografia row (row 22)
# - Fixes the column A/B <col> width overlap and several row heights

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Fix the overlapping column definition (col A was min=1,max=2; should
#    be min=1,max=1 so column B's own 60.71 width/style definition wins).
#    Touching column B's width causes the engine to split the ranges,
#    which also fixes column A's span.
# ---------------------------------------------------------------------
$ws.Columns.Item(2).ColumnWidth = 60.7109375

# ---------------------------------------------------------------------
# 2. Row 10 (Objetivos:) - B/C get the new Portuguese objectives text
#    (previously held the "198273 - Domingos Savio Giordani" text, which
#    moves down to row 13).
# ---------------------------------------------------------------------
$ws.Range("B10").Value = "Conduzir os alunos no desenvolvimento de um projeto de conclusão de curso sobre tema específico relacionado a engenharia de produção."
$ws.Range("C10").Value = "Conduzir os alunos no desenvolvimento de um projeto de conclusão de curso sobre tema específico relacionado a engenharia de produção."

# ---------------------------------------------------------------------
# 3. Row 13 becomes the "Docentes responsáveis:" value row (no A label
#    cell anymore - the label text now lives on row 12 (A12, unchanged)).
# ---------------------------------------------------------------------
$ws.Range("A13").ClearContents()
$ws.Range("B13").Value = "198273 - Domingos Savio Giordani"
$ws.Range("C13").Value = "198273 - Domingos Savio Giordani"
$ws.Rows.Item(13).AutoFit()

# ---------------------------------------------------------------------
# 4. Row 14 becomes "Programa resumido:" with the new PT short syllabus.
# ---------------------------------------------------------------------
$ws.Range("A14").Value = "Programa resumido:"
$ws.Range("B14").Value = "1) Metodologia Cientifica. 2) Projeto de Monografia. 3) Métodos de Pesquisa. 4) Normas de um Projeto de Pesquisa."
$ws.Range("C14").Value = "1) Metodologia Cientifica. 2) Projeto de Monografia. 3) Métodos de Pesquisa. 4) Normas de um Projeto de Pesquisa."
$ws.Rows.Item(14).RowHeight = 60

# ---------------------------------------------------------------------
# 5. Row 15 becomes "Short syllabus:" (English) - moved down from row 14.
# ---------------------------------------------------------------------
$ws.Range("A15").Value = "Short syllabus:"
$ws.Range("B15").Value = "1) Scientific Methodology. 2) Monograph Project. 3) Research Methods. 4) Norms of a Research Project."
$ws.Range("C15").Value = "1) Scientific Methodology. 2) Monograph Project. 3) Research Methods. 4) Norms of a Research Project."
$ws.Rows.Item(15).RowHeight = 60

# ---------------------------------------------------------------------
# 6. Row 16 becomes "Programa:" with the new PT full syllabus text
#    (row height stays 120).
# ---------------------------------------------------------------------
$ws.Range("A16").Value = "Programa:"
$ws.Range("B16").Value = "1  Metodologia Cientifica: Concepção e definição. 2  Monografia Cientifica: O que é um projeto de pesquisa. As etapas de um projeto de pesquisa. 3  Métodos de Pesquisa utilizados na Engenharia de Produção  4- Pontos essenciais de um projeto de TCC. 5  Pontos essenciais de uma monografia de TCC. 6  Normas para elaboração de Referências Bibliográficas."
$ws.Range("C16").Value = "1  Metodologia Cientifica: Concepção e definição. 2  Monografia Cientifica: O que é um projeto de pesquisa. As etapas de um projeto de pesquisa. 3  Métodos de Pesquisa utilizados na Engenharia de Produção  4- Pontos essenciais de um projeto de TCC. 5  Pontos essenciais de uma monografia de TCC. 6  Normas para elaboração de Referências Bibliográficas."
$ws.Rows.Item(16).RowHeight = 120

# ---------------------------------------------------------------------
# 7. Row 17 becomes "Syllabus:" (English) with new B/C cells (previously
#    row 17 only had the A-label, no B/C).
# ---------------------------------------------------------------------
$ws.Range("A17").Value = "Syllabus:"
$ws.Range("B16").Copy() | Out-Null
$ws.Range("B17").PasteSpecial(-4122) | Out-Null
$ws.Range("C16").Copy() | Out-Null
$ws.Range("C17").PasteSpecial(-4122) | Out-Null
$ws.Range("B17").Value = "1 - Scientific Methodology: Conception and definition. 2 - Scientific Monograph: What is a research project. The steps of a research project. 3 - Research Methods used in Industrial Engineering 4- Essential points of a TCC project. 5 - Essential points of a TCC monograph. 6 - Norms for elaboration of Bibliographical References."
$ws.Range("C17").Value = "1 - Scientific Methodology: Conception and definition. 2 - Scientific Monograph: What is a research project. The steps of a research project. 3 - Research Methods used in Industrial Engineering 4- Essential points of a TCC project. 5 - Essential points of a TCC monograph. 6 - Norms for elaboration of Bibliographical References."
$ws.Rows.Item(17).RowHeight = 120

# ---------------------------------------------------------------------
# 8. Row 18 becomes "Avaliação:" label only (B/C cleared - that value
#    moved to row 19 as "Método:").
# ---------------------------------------------------------------------
$ws.Range("A18").Value = "Avaliação:"
$ws.Range("B18").ClearContents()
$ws.Range("C18").ClearContents()
$ws.Rows.Item(18).AutoFit()

# ---------------------------------------------------------------------
# 9. Row 19 becomes "Método:" with the Desenvolvimento... text.
# ---------------------------------------------------------------------
$ws.Range("A19").Value = "Método:"
$ws.Range("B19").Value = "Desenvolvimento e apresentação do Projeto monografia a ser desenvolvida na disciplina de Trabalho de Graduação em Engenharia de Produção II, conforme norma do Curso de Engenharia de Produção"
$ws.Range("C19").Value = "Desenvolvimento e apresentação do Projeto monografia a ser desenvolvida na disciplina de Trabalho de Graduação em Engenharia de Produção II, conforme norma do Curso de Engenharia de Produção"
$ws.Rows.Item(19).RowHeight = 60

# ---------------------------------------------------------------------
# 10. Row 20 becomes "Critério:" with the Avaliação Ad hoc text.
# ---------------------------------------------------------------------
$ws.Range("A20").Value = "Critério:"
$ws.Range("B20").Value = "Avaliação Ad hoc por 2 examinadores. A nota da disciplina será a média das duas notas"
$ws.Range("C20").Value = "Avaliação Ad hoc por 2 examinadores. A nota da disciplina será a média das duas notas"
$ws.Rows.Item(20).RowHeight = 60

# ---------------------------------------------------------------------
# 11. Row 21 becomes "Norma de recuperação:" with the Reapresentação
#     text (row height now 60, was 120).
# ---------------------------------------------------------------------
$ws.Range("A21").Value = "Norma de recuperação:"
$ws.Range("B21").Value = "Reapresentação do trabalho modificado para nova avaliação"
$ws.Range("C21").Value = "Reapresentação do trabalho modificado para nova avaliação"
$ws.Rows.Item(21).RowHeight = 60

# ---------------------------------------------------------------------
# 12. New row 22: "Bibliografia:" with the bibliography text (120pt tall).
# ---------------------------------------------------------------------
$ws.Range("A21").Copy() | Out-Null
$ws.Range("A22").PasteSpecial(-4122) | Out-Null
$ws.Range("B21").Copy() | Out-Null
$ws.Range("B22").PasteSpecial(-4122) | Out-Null
$ws.Range("C21").Copy() | Out-Null
$ws.Range("C22").PasteSpecial(-4122) | Out-Null
$ws.Range("A22").Value = "Bibliografia:"
$ws.Range("B22").Value = "Cauchick-Miguel, P. A.,   Metodologia de pesquisa em engenharia de produção e gestão de operações / Afonso Fleury ... [et al.] ; coordenação . - 3. ed. - Rio de Janeiro : Elsevier, 2018. Cauchick-Miguel, P. A.,   Metodologia de pesquisa em engenharia , 1. ed. - Rio de Janeiro : GEN LTC, 2019. BOOTH, W.; COLOMB, G.; WILLIAMS, J. A arte da Pesquisa. 3 ed. Martins Fontes. São Paulo. 2005.GIL, A.C. Como elaborar projetos de pesquisa. 5 ed. Atlas, São Paulo, 2010.MEDEIROS, J. B. Redação Cientifica: A Prática de Fichamentos, Resumos e Resenhas. 11 ed. São Paulo: Atlas, 2009"
$ws.Range("C22").Value = "Cauchick-Miguel, P. A.,   Metodologia de pesquisa em engenharia de produção e gestão de operações / Afonso Fleury ... [et al.] ; coordenação . - 3. ed. - Rio de Janeiro : Elsevier, 2018. Cauchick-Miguel, P. A.,   Metodologia de pesquisa em engenharia , 1. ed. - Rio de Janeiro : GEN LTC, 2019. BOOTH, W.; COLOMB, G.; WILLIAMS, J. A arte da Pesquisa. 3 ed. Martins Fontes. São Paulo. 2005.GIL, A.C. Como elaborar projetos de pesquisa. 5 ed. Atlas, São Paulo, 2010.MEDEIROS, J. B. Redação Cientifica: A Prática de Fichamentos, Resumos e Resenhas. 11 ed. São Paulo: Atlas, 2009"
$ws.Rows.Item(22).RowHeight = 120

Write-Output "Done"
